$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9956793189048767
$ws.Range("B1").Value = 2.736323595046997
$ws.Range("C1").Value = 4.659855365753174
$ws.Range("D1").Value = 1.094549536705017
$ws.Range("E1").Value = 1.27775502204895
